$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2456.4285
$ws.Range("J17").Value = 2456.4285
$ws.Range("L17").Value = 7369.2855
$ws.Range("N17").Value = -7705.2855
# Row 51
$ws.Range("H51").Value = 2663.6365
$ws.Range("I51").Value = 2642.8572
$ws.Range("J51").Value = 2700
$ws.Range("K51").Value = 2642.8572
$ws.Range("L51").Value = 2700
$ws.Range("M51").Value = -2158.8572
$ws.Range("N51").Value = -3668
# Row 86
$ws.Range("H86").Value = 2779.2
$ws.Range("I86").Value = 2989.9473
$ws.Range("K86").Value = 2989.9473
$ws.Range("M86").Value = -1866.9473
# Row 89
$ws.Range("H89").Value = 2779.2
$ws.Range("I89").Value = 2989.9473
$ws.Range("K89").Value = 14949.7365
$ws.Range("M89").Value = -9333.736499999999
# Row 127
$ws.Range("H127").Value = 1355.4546
$ws.Range("I127").Value = 1335.6666
$ws.Range("K127").Value = 4006.9998
$ws.Range("M127").Value = 953.0001999999999
# Row 138
$ws.Range("H138").Value = 2437.283
$ws.Range("I138").Value = 1354.5416
$ws.Range("J138").Value = 3333.3447
$ws.Range("K138").Value = 4063.6248
$ws.Range("L138").Value = 10000.0341
$ws.Range("M138").Value = 1076.3752
$ws.Range("N138").Value = -20280.0341

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 343.0909
$ws.Range("I5").Value = 287.25
$ws.Range("K5").Value = 287.25
$ws.Range("M5").Value = -175.25
# Row 32
$ws.Range("H32").Value = 2986.587
$ws.Range("I32").Value = 2849.3076
$ws.Range("K32").Value = 2849.3076
$ws.Range("M32").Value = -2562.3076
# Row 37
$ws.Range("H37").Value = 34726.855
$ws.Range("J37").Value = 34726.855
$ws.Range("L37").Value = 34726.855
$ws.Range("N37").Value = -35272.855
# Row 41
$ws.Range("H41").Value = 7126.25
$ws.Range("I41").Value = 7126.25
$ws.Range("K41").Value = 7126.25
$ws.Range("M41").Value = -6712.25
# Row 45
$ws.Range("H45").Value = 1503.7778
$ws.Range("I45").Value = 2415.25
$ws.Range("K45").Value = 2415.25
$ws.Range("M45").Value = -2038.25
# Row 63
$ws.Range("H63").Value = 4750.3477
$ws.Range("I63").Value = 4356.3335
$ws.Range("J63").Value = 5489.125
$ws.Range("K63").Value = 4356.3335
$ws.Range("L63").Value = 5489.125
$ws.Range("M63").Value = -3670.3335
$ws.Range("N63").Value = -6861.125
# Row 66
$ws.Range("H66").Value = 4750.3477
$ws.Range("I66").Value = 4356.3335
$ws.Range("J66").Value = 5489.125
$ws.Range("K66").Value = 21781.6675
$ws.Range("L66").Value = 27445.625
$ws.Range("M66").Value = -18349.6675
$ws.Range("N66").Value = -34309.625
# Row 74
$ws.Range("H74").Value = 1733.32
$ws.Range("I74").Value = 1540.381
$ws.Range("J74").Value = 2746.25
$ws.Range("K74").Value = 1540.381
$ws.Range("L74").Value = 2746.25
$ws.Range("M74").Value = -666.3810000000001
$ws.Range("N74").Value = -4494.25
# Row 77
$ws.Range("H77").Value = 1733.32
$ws.Range("I77").Value = 1540.381
$ws.Range("J77").Value = 2746.25
$ws.Range("K77").Value = 7701.905000000001
$ws.Range("L77").Value = 13731.25
$ws.Range("M77").Value = -3333.905000000001
$ws.Range("N77").Value = -22467.25
# Row 88
$ws.Range("H88").Value = 3165
$ws.Range("I88").Value = 2502.5
$ws.Range("K88").Value = 2502.5
$ws.Range("M88").Value = -2096.5
# Row 91
$ws.Range("H91").Value = 3165
$ws.Range("I91").Value = 2502.5
$ws.Range("K91").Value = 2502.5
$ws.Range("M91").Value = -1098.5
# Row 110
$ws.Range("H110").Value = 1978.4166
$ws.Range("I110").Value = 705.5
$ws.Range("K110").Value = 705.5
$ws.Range("M110").Value = 1339.5
# Row 132
$ws.Range("H132").Value = 1932
$ws.Range("I132").Value = 1952.174
$ws.Range("K132").Value = 5856.522
$ws.Range("M132").Value = -3326.522

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 343.0909
$ws.Range("I4").Value = 287.25
$ws.Range("K4").Value = 287.25
$ws.Range("M4").Value = -172.25
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 16
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("K16").Value = 500
$ws.Range("M16").Value = -330
# Row 22
$ws.Range("H22").Value = 1856.4286
$ws.Range("I22").Value = 1832.1666
$ws.Range("K22").Value = 1832.1666
$ws.Range("M22").Value = -1659.1666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 26
$ws.Range("H26").Value = 11210
$ws.Range("J26").Value = 13625
$ws.Range("L26").Value = 13625
$ws.Range("N26").Value = -14199
# Row 56
$ws.Range("H56").Value = 235.75
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
# Row 58
$ws.Range("H58").Value = 2205.68
$ws.Range("I58").Value = 2293.8125
$ws.Range("J58").Value = 2049
$ws.Range("K58").Value = 2293.8125
$ws.Range("L58").Value = 2049
$ws.Range("M58").Value = -2090.8125
$ws.Range("N58").Value = -2455
# Row 136
$ws.Range("H136").Value = 2205.68
$ws.Range("I136").Value = 2293.8125
$ws.Range("J136").Value = 2049
$ws.Range("K136").Value = 6881.4375
$ws.Range("L136").Value = 6147
$ws.Range("M136").Value = -4331.4375
$ws.Range("N136").Value = -11247

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1343.5294
$ws.Range("I107").Value = 2032.8334
$ws.Range("J107").Value = 967.5454999999999
$ws.Range("K107").Value = 6098.5002
$ws.Range("L107").Value = 2902.6365
$ws.Range("M107").Value = -4178.5002
$ws.Range("N107").Value = -6742.6365
# Row 132
$ws.Range("H132").Value = 2033.3334
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 22500
$ws.Range("N132").Value = -27560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 23771
$ws.Range("J15").Value = 23771
$ws.Range("L15").Value = 23771
$ws.Range("N15").Value = -24347
# Row 23
$ws.Range("H23").Value = 9450
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10446
# Row 80
$ws.Range("H80").Value = 6443.273
$ws.Range("J80").Value = 7333.3335
$ws.Range("L80").Value = 7333.3335
$ws.Range("N80").Value = -9329.333500000001
# Row 81
$ws.Range("H81").Value = 23771
$ws.Range("J81").Value = 23771
$ws.Range("L81").Value = 23771
$ws.Range("N81").Value = -25767
# Row 83
$ws.Range("H83").Value = 6443.273
$ws.Range("J83").Value = 7333.3335
$ws.Range("L83").Value = 36666.6675
$ws.Range("N83").Value = -46650.6675
# Row 84
$ws.Range("H84").Value = 23771
$ws.Range("J84").Value = 23771
$ws.Range("L84").Value = 71313
$ws.Range("N84").Value = -81297

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 366.77777
$ws.Range("I16").Value = 350.125
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 350.125
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -180.125
$ws.Range("N16").Value = -840
# Row 55
$ws.Range("H55").Value = 1672.3334
$ws.Range("I55").Value = 224.45454
$ws.Range("J55").Value = 3947.5715
$ws.Range("K55").Value = 224.45454
$ws.Range("L55").Value = 3947.5715
$ws.Range("M55").Value = -51.45454000000001
$ws.Range("N55").Value = -4293.5715
# Row 93
$ws.Range("H93").Value = 3876.0715
$ws.Range("I93").Value = 830.7273
$ws.Range("K93").Value = 830.7273
$ws.Range("M93").Value = 417.2727

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 14240.172
$ws.Range("J62").Value = 14367.883
$ws.Range("L62").Value = 14367.883
$ws.Range("N62").Value = -15615.883
# Row 65
$ws.Range("H65").Value = 14240.172
$ws.Range("J65").Value = 14367.883
$ws.Range("L65").Value = 71839.41499999999
$ws.Range("N65").Value = -78079.41499999999
# Row 75
$ws.Range("H75").Value = 500043300
$ws.Range("J75").Value = 500043300
$ws.Range("L75").Value = 500043300
$ws.Range("N75").Value = -500045172
# Row 78
$ws.Range("H78").Value = 500043300
$ws.Range("J78").Value = 500043300
$ws.Range("L78").Value = 1500129900
$ws.Range("N78").Value = -1500139260
# Row 81
$ws.Range("H81").Value = 8200.333000000001
$ws.Range("I81").Value = 4601
$ws.Range("K81").Value = 9202
$ws.Range("M81").Value = -8141
# Row 84
$ws.Range("H84").Value = 8200.333000000001
$ws.Range("I84").Value = 4601
$ws.Range("K84").Value = 46010
$ws.Range("M84").Value = -40706
# Row 136
$ws.Range("H136").Value = 1308.2559
$ws.Range("I136").Value = 913.5405
$ws.Range("J136").Value = 3742.3333
$ws.Range("K136").Value = 2740.6215
$ws.Range("L136").Value = 11226.9999
$ws.Range("M136").Value = -190.6214999999997
$ws.Range("N136").Value = -16326.9999
